# Ajuste para permitir parametrizar horarios de los días viernes
#
# The sheet currently has a single "L-V" (Mon-Fri) entry/exit pair in
# columns I/J, followed by the Saturday entry/exit pair in columns K/L.
# We split the Mon-Fri pair into a Mon-Thu pair (renamed, keeping its
# current columns I/J) and a new, separate Friday pair (new columns K/L),
# pushing the existing Saturday pair two columns to the right (M/N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the width currently used by column J so the two freshly
# inserted columns (which will hold the new Friday headers) match it.
$existingWidth = $ws.Range("J1").ColumnWidth

# Insert two new blank columns at K:L -- this shifts the current
# Saturday entry/exit headers (K1:L1) two columns to the right, to M1:N1,
# and leaves K1:L1 empty and ready for the new Friday headers.
$ws.Range("K1:L1").EntireColumn.Insert()

# Mon-Fri (L-V) becomes Mon-Thu (L-J).
$ws.Range("I1").Value = "HORA DE ENTRADA L-J"
$ws.Range("J1").Value = "HORA DE SALIDA L-J"

# New, dedicated Friday schedule columns.
$ws.Range("K1").Value = "HORA DE ENTRADA VIERNES"
$ws.Range("L1").Value = "HORA DE SALIDA VIERNES"

# Match the new columns' width to the neighboring column so the sheet
# keeps a consistent look.
$ws.Range("K1:L1").ColumnWidth = $existingWidth

# Leave the same selection state the saved workbook shipped with.
[void]$ws.Range("E7").Select()
